# Update the "想去人数" (interested-people count) figures in column F
# for the sheets that list exhibition data: "展览" and "全部类型".
# Both sheets mirror the same rows, so the same row -> new value map applies.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 197
    5  = 3422
    6  = 5968
    7  = 326
    11 = 8794
    12 = 2354
    13 = 248
    14 = 5352
    15 = 10255
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
